$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.440.42"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "1.871.91"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8920"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07967"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.93%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "1.845.55"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.453"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.638"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "93.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008972"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "27.464.02"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.191"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "2.080.17"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.881"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.110"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.183"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08933"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7590"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.85%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.164"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.523"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05310"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01965"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.229"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5258"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1651"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.665"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.24%  "
